$wb = $excel.ActiveWorkbook

# --- Fix the FM mkdocs table: it was skipping the "Storage node file" row
# that already exists in the Source table (row 41), which hid it from the
# published FM docs table. Insert a new row at position 33 and fill it with
# the same formula pattern as its neighbours, pointing at Source table row 41.
$ws3 = $wb.Worksheets.Item("FM mkdocs table")

$ws3.Rows.Item(33).Insert()

$dq = [char]34

$ws3.Range("A33").Formula = "=IF(ISBLANK('Source table'!A41)," + $dq + " " + $dq + ",'Source table'!A41)"
$ws3.Range("B33").Formula = "=IFERROR(VLOOKUP('Source table'!B41,'mkdocs symbols'!`$A`$1:`$C`$5,2,0)," + $dq + " " + $dq + ")"
$ws3.Range("C33").Formula = "=IFERROR(VLOOKUP('Source table'!C41,'mkdocs symbols'!`$A`$1:`$C`$5,2,0)," + $dq + " " + $dq + ")"
$ws3.Range("D33").Formula = "=IF(ISBLANK('Source table'!D41)," + $dq + " " + $dq + ",'Source table'!D41)"
$ws3.Range("E33").Formula = "=IF(OR(ISBLANK('Source table'!E41),ISBLANK('Source table'!F41))," + $dq + " " + $dq + "," + $dq + "[" + $dq + "&'Source table'!F41&" + $dq + "][" + $dq + "&'Source table'!E41&" + $dq + "." + $dq + "&'Source table'!F41&" + $dq + "]" + $dq + ")"
$ws3.Range("F33").Formula = "=IF(ISBLANK('Source table'!G41)," + $dq + " " + $dq + "," + $dq + "_" + $dq + "&'Source table'!G41&" + $dq + "_" + $dq + ")"

# (Row.Insert() already carried the neighbouring rows' styles onto row 33,
# same as Excel does, so B33/C33 stay unstyled and A33/D33/E33/F33 keep
# style index 3 - no extra formatting work needed.)

# --- Restore view state: select the source row that was being hidden on the
# "Source table" sheet, then leave "FM mkdocs table" as the active sheet/tab.
$ws4 = $wb.Worksheets.Item("Source table")
$ws4.Activate()
$ws4.Range("A41").Select()

$ws3.Activate()
$ws3.Range("A51").Select()
